# Auto-generated Excel COM-interop script to update Typhon_Profits price data
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 215.41667
$ws.Range("I9").Value = 77.666664
$ws.Range("J9").Value = 353.16666
$ws.Range("K9").Value = 77.666664
$ws.Range("L9").Value = 353.16666
$ws.Range("M9").Value = 91.333336
$ws.Range("N9").Value = -691.16666
$ws.Range("H15").Value = 1564.2988
$ws.Range("I15").Value = 1564.2988
$ws.Range("K15").Value = 4692.8964
$ws.Range("M15").Value = -4523.8964
$ws.Range("H32").Value = 503.33334
$ws.Range("I32").Value = 562.5
$ws.Range("J32").Value = 456
$ws.Range("K32").Value = 562.5
$ws.Range("L32").Value = 456
$ws.Range("M32").Value = -236.5
$ws.Range("N32").Value = -1108
$ws.Range("H92").Value = 90910090
$ws.Range("J92").Value = 774.5
$ws.Range("L92").Value = 774.5
$ws.Range("N92").Value = -3270.5
$ws.Range("H96").Value = 25000352
$ws.Range("I96").Value = 25000352
$ws.Range("K96").Value = 75001056
$ws.Range("M96").Value = -74999683
$ws.Range("H111").Value = 3421.8333
$ws.Range("J111").Value = 1482.75
$ws.Range("L111").Value = 4448.25
$ws.Range("N111").Value = -10582.25
$ws.Range("H137").Value = 28078.395
$ws.Range("I137").Value = 1616.0416
$ws.Range("J137").Value = 73442.42999999999
$ws.Range("K137").Value = 4848.1248
$ws.Range("L137").Value = 220327.29
$ws.Range("M137").Value = -2298.1248
$ws.Range("N137").Value = -225427.29
$ws.Range("H138").Value = 2385.2239
$ws.Range("J138").Value = 2284.07
$ws.Range("L138").Value = 6852.210000000001
$ws.Range("N138").Value = -17132.21
$ws.Range("H141").Value = 1293.6
$ws.Range("I141").Value = 951.32355
$ws.Range("J141").Value = 3233.1667
$ws.Range("K141").Value = 2853.97065
$ws.Range("L141").Value = 9699.500100000001
$ws.Range("M141").Value = 2326.02935
$ws.Range("N141").Value = -20059.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2800
$ws.Range("J3").Value = 2800
$ws.Range("L3").Value = 2800
$ws.Range("N3").Value = -3030
$ws.Range("H45").Value = 2443.838
$ws.Range("I45").Value = 2179.5862
$ws.Range("J45").Value = 3401.75
$ws.Range("K45").Value = 2179.5862
$ws.Range("L45").Value = 3401.75
$ws.Range("M45").Value = -1802.5862
$ws.Range("N45").Value = -4155.75
$ws.Range("H74").Value = 58824520
$ws.Range("I74").Value = 90910080
$ws.Range("K74").Value = 90910080
$ws.Range("M74").Value = -90909206
$ws.Range("H77").Value = 58824520
$ws.Range("I77").Value = 90910080
$ws.Range("K77").Value = 454550400
$ws.Range("M77").Value = -454546032
$ws.Range("H98").Value = 17250
$ws.Range("J98").Value = 17250
$ws.Range("L98").Value = 17250
$ws.Range("N98").Value = -23240
$ws.Range("H132").Value = 29174.8
$ws.Range("I132").Value = 1645.9783
$ws.Range("K132").Value = 4937.9349
$ws.Range("M132").Value = -2407.9349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1345.0588
$ws.Range("I99").Value = 1488.5714
$ws.Range("J99").Value = 1244.6
$ws.Range("K99").Value = 1488.5714
$ws.Range("L99").Value = 1244.6
$ws.Range("M99").Value = 9.42859999999996
$ws.Range("N99").Value = -4240.6
$ws.Range("H105").Value = 2546.4546
$ws.Range("I105").Value = 2454.8
$ws.Range("J105").Value = 2742.8572
$ws.Range("K105").Value = 2454.8
$ws.Range("L105").Value = 2742.8572
$ws.Range("M105").Value = -707.8000000000002
$ws.Range("N105").Value = -6236.8572
$ws.Range("H107").Value = 903.12
$ws.Range("I107").Value = 780.26666
$ws.Range("J107").Value = 1087.4
$ws.Range("K107").Value = 780.26666
$ws.Range("L107").Value = 1087.4
$ws.Range("M107").Value = 1139.73334
$ws.Range("N107").Value = -4927.4
$ws.Range("H134").Value = 21221.629
$ws.Range("I134").Value = 24547.479
$ws.Range("K134").Value = 73642.43700000001
$ws.Range("M134").Value = -71107.43700000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("H62").Value = 333338340
$ws.Range("I62").Value = 333338340
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 333338340
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = -333337716
$ws.Range("H65").Value = 333338340
$ws.Range("I65").Value = 333338340
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 1666691700
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = -1666688580
$ws.Range("H105").Value = 8929368
$ws.Range("I105").Value = 9616165
$ws.Range("J105").Value = 999
$ws.Range("K105").Value = 9616165
$ws.Range("L105").Value = 999
$ws.Range("M105").Value = -9614418
$ws.Range("N105").Value = -4493
$ws.Range("H132").Value = 13428.272
$ws.Range("I132").Value = 14545.263
$ws.Range("K132").Value = 43635.789
$ws.Range("M132").Value = -41105.789

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 10225
$ws.Range("I107").Value = 100000
$ws.Range("J107").Value = 250
$ws.Range("K107").Value = 300000
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = -298080
$ws.Range("N107").Value = -4590
$ws.Range("H109").Value = 4590.364
$ws.Range("I109").Value = 1400
$ws.Range("J109").Value = 6413.4287
$ws.Range("K109").Value = 4200
$ws.Range("L109").Value = 19240.2861
$ws.Range("M109").Value = -3160
$ws.Range("N109").Value = -21320.2861
$ws.Range("H131").Value = 766.8099999999999
$ws.Range("I131").Value = 338.33334
$ws.Range("J131").Value = 794.15955
$ws.Range("K131").Value = 1015.00002
$ws.Range("L131").Value = 2382.47865
$ws.Range("M131").Value = 4024.99998
$ws.Range("N131").Value = -12462.47865

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10378.571
$ws.Range("H73").Value = 10378.571
$ws.Range("H97").Value = 915
$ws.Range("H102").Value = 17859226
$ws.Range("I102").Value = 23811828
$ws.Range("J102").Value = 1416.2858
$ws.Range("K102").Value = 23811828
$ws.Range("L102").Value = 1416.2858
$ws.Range("M102").Value = -23810206
$ws.Range("N102").Value = -4660.2858
$ws.Range("H122").Value = 70176376
$ws.Range("I122").Value = 22222842
$ws.Range("J122").Value = 250002130
$ws.Range("K122").Value = 66668526
$ws.Range("L122").Value = 750006390
$ws.Range("M122").Value = -66666076
$ws.Range("N122").Value = -750011290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 799.8
$ws.Range("I16").Value = 799.8
$ws.Range("K16").Value = 799.8
$ws.Range("M16").Value = -629.8
$ws.Range("H61").Value = 5923.636
$ws.Range("I61").Value = 2622.8572
$ws.Range("K61").Value = 2622.8572
$ws.Range("M61").Value = -2420.8572
$ws.Range("H100").Value = 1899
$ws.Range("I100").Value = 1360
$ws.Range("J100").Value = 2284
$ws.Range("K100").Value = 1360
$ws.Range("L100").Value = 2284
$ws.Range("M100").Value = -819
$ws.Range("N100").Value = -3366
$ws.Range("H113").Value = 5923.636
$ws.Range("I113").Value = 2622.8572
$ws.Range("K113").Value = 2622.8572
$ws.Range("M113").Value = -452.8571999999999
$ws.Range("H132").Value = 3907.4285
$ws.Range("I132").Value = 3088.5
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 9265.5
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -6735.5
$ws.Range("N132").Value = -20058.0005
$ws.Range("H136").Value = 35030.066
$ws.Range("I136").Value = 42953.418
$ws.Range("J136").Value = 3336.6667
$ws.Range("K136").Value = 128860.254
$ws.Range("L136").Value = 10010.0001
$ws.Range("M136").Value = -126310.254
$ws.Range("N136").Value = -15110.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1426.8
$ws.Range("I132").Value = 1040.1875
$ws.Range("J132").Value = 2973.25
$ws.Range("K132").Value = 3120.5625
$ws.Range("L132").Value = 8919.75
$ws.Range("M132").Value = -590.5625
$ws.Range("N132").Value = -13979.75
$ws.Range("H136").Value = 32259662
$ws.Range("I136").Value = 37038468
$ws.Range("J136").Value = 2725
$ws.Range("K136").Value = 111115404
$ws.Range("L136").Value = 8175
$ws.Range("M136").Value = -111112854
$ws.Range("N136").Value = -13275

Write-Output "Applied 208 value updates across 8 sheets"